# DataEngine.xlsx - 1st Commit on 3rd May
# Adds PASS/FAIL result columns to the TestCases and TestSteps sheets,
# and leaves the selection/active-tab on the TestSteps sheet.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestSteps = $wb.Worksheets.Item("TestSteps")

# --- TestCases sheet: column D gets a PASS/PASS/FAIL result ---
$wsTestCases.Range("D2").Value = "PASS"
$wsTestCases.Range("D3").Value = "PASS"
$wsTestCases.Range("D4").Value = "FAIL"

# --- TestSteps sheet: column H gets PASS for every data row (2-33) ---
for ($row = 2; $row -le 33; $row++) {
    $wsTestSteps.Cells.Item($row, 8).Value = "PASS"
}

# --- View state: TestCases selection moves to D2:D4, no longer the active tab ---
$wsTestCases.Activate()
$wsTestCases.Range("D2:D4").Select()

# --- View state: TestSteps becomes the active tab, selection moves to G35 ---
$wsTestSteps.Activate()
$wsTestSteps.Range("G35").Select()
